{"js": "// Replace the 25 two-digit-division problems in the document body.\n// Each entry maps the original \"a\u00f7b=\" text to its replacement.\nconst replacements = [\n  [\"79\u00f72=\", \"33\u00f74=\"],\n  [\"79\u00f78=\", \"18\u00f76=\"],\n  [\"19\u00f76=\", \"60\u00f72=\"],\n  [\"40\u00f78=\", \"73\u00f75=\"],\n  [\"35\u00f79=\", \"87\u00f72=\"],\n  [\"85\u00f72=\", \"78\u00f73=\"],\n  [\"71\u00f75=\", \"88\u00f77=\"],\n  [\"50\u00f79=\", \"62\u00f73=\"],\n  [\"13\u00f74=\", \"57\u00f74=\"],\n  [\"98\u00f72=\", \"55\u00f74=\"],\n  [\"60\u00f75=\", \"27\u00f76=\"],\n  [\"26\u00f75=\", \"92\u00f78=\"],\n  [\"38\u00f74=\", \"71\u00f78=\"],\n  [\"43\u00f74=\", \"57\u00f73=\"],\n  [\"15\u00f73=\", \"65\u00f76=\"],\n  [\"39\u00f73=\", \"35\u00f78=\"],\n  [\"63\u00f79=\", \"47\u00f72=\"],\n  [\"21\u00f78=\", \"18\u00f75=\"],\n  [\"33\u00f79=\", \"39\u00f78=\"],\n  [\"33\u00f75=\", \"35\u00f74=\"],\n  [\"34\u00f79=\", \"10\u00f77=\"],\n  [\"81\u00f79=\", \"69\u00f74=\"],\n  [\"76\u00f75=\", \"53\u00f74=\"],\n  [\"80\u00f77=\", \"40\u00f75=\"],\n  [\"94\u00f78=\", \"86\u00f73=\"],\n];\n\nfor (const [before, after] of replacements) {\n  const results = context.document.body.search(before, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(after, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the 25 two-digit-division problems in the document body.\n# Each entry maps the original \"a\u00f7b=\" text to its replacement.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"79\u00f72=\", \"33\u00f74=\"),\n    @(\"79\u00f78=\", \"18\u00f76=\"),\n    @(\"19\u00f76=\", \"60\u00f72=\"),\n    @(\"40\u00f78=\", \"73\u00f75=\"),\n    @(\"35\u00f79=\", \"87\u00f72=\"),\n    @(\"85\u00f72=\", \"78\u00f73=\"),\n    @(\"71\u00f75=\", \"88\u00f77=\"),\n    @(\"50\u00f79=\", \"62\u00f73=\"),\n    @(\"13\u00f74=\", \"57\u00f74=\"),\n    @(\"98\u00f72=\", \"55\u00f74=\"),\n    @(\"60\u00f75=\", \"27\u00f76=\"),\n    @(\"26\u00f75=\", \"92\u00f78=\"),\n    @(\"38\u00f74=\", \"71\u00f78=\"),\n    @(\"43\u00f74=\", \"57\u00f73=\"),\n    @(\"15\u00f73=\", \"65\u00f76=\"),\n    @(\"39\u00f73=\", \"35\u00f78=\"),\n    @(\"63\u00f79=\", \"47\u00f72=\"),\n    @(\"21\u00f78=\", \"18\u00f75=\"),\n    @(\"33\u00f79=\", \"39\u00f78=\"),\n    @(\"33\u00f75=\", \"35\u00f74=\"),\n    @(\"34\u00f79=\", \"10\u00f77=\"),\n    @(\"81\u00f79=\", \"69\u00f74=\"),\n    @(\"76\u00f75=\", \"53\u00f74=\"),\n    @(\"80\u00f77=\", \"40\u00f75=\"),\n    @(\"94\u00f78=\", \"86\u00f73=\")\n)\n\nforeach ($pair in $replacements) {\n    $before = $pair[0]\n    $after = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $before\n    $find.Replacement.Text = $after\n    $find.Execute($before, $false, $true, $false, $false, $false, $true, 1, $false, $after, 2)\n}\n"}
